$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.959.56"
$ws.Range("E2").Value = "  -3.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.719.09"
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.25"
$ws.Range("E5").Value = "  -6.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4848"
$ws.Range("E7").Value = "  +6.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3491"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.13"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07249"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.87"
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.850"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.722.96"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.851"
$ws.Range("E16").Value = "  -4.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.65"
$ws.Range("E17").Value = "  -6.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06373"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.50"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.637"
$ws.Range("E22").Value = "  -2.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.005.93"
$ws.Range("E23").Value = "  -3.32%  "
$ws.Range("E24").Value = "  -3.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.082"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.67"
$ws.Range("E26").Value = "  -5.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.92"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.913.62"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.065"
$ws.Range("E29").Value = "  -4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.94"
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.028"
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09272"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("E33").Value = "  -1.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.341"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05883"
$ws.Range("E35").Value = "  -3.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02173"
$ws.Range("E36").Value = "  -4.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.445"
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.97"
$ws.Range("E38").Value = "  -7.03%  "
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.734"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("B41").Value = "Frax"
$ws.Range("C41").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.0000"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5996"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.087"
$ws.Range("E43").Value = "  -7.99%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.528"
$ws.Range("E44").Value = "  -4.31%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.69"
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.572"
$ws.Range("E46").Value = "  -4.34%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5617"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "117.70"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.832"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.107"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06648"
$ws.Range("E51").Value = "  -2.34%  "
